{"js": "const oldText = \"Datas da campa\u00f1a de Constelaci\u00f3n de H\u00e9rcules 2022: 13-22 de xu\u00f1o, 12-21 de xullo, 10-19 de agosto\";\nconst newText = \"Datas da campa\u00f1a de 2022 que usan Constelaci\u00f3n de H\u00e9rcules: 13-22 de xu\u00f1o, 12-21 de xullo, 10-19 de agosto\";\n\nconst results = context.document.body.search(oldText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$oldText = \"Datas da campa\u00f1a de Constelaci\u00f3n de H\u00e9rcules 2022: 13-22 de xu\u00f1o, 12-21 de xullo, 10-19 de agosto\"\n$newText = \"Datas da campa\u00f1a de 2022 que usan Constelaci\u00f3n de H\u00e9rcules: 13-22 de xu\u00f1o, 12-21 de xullo, 10-19 de agosto\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $oldText\n$find.Replacement.Text = $newText\n$find.Forward = $true\n$find.Wrap = 2\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$find.Execute(\n    $oldText,\n    $false,\n    $true,\n    $false,\n    $false,\n    $false,\n    $true,\n    2,\n    $false,\n    $newText,\n    2\n)\n"}
